$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8742.895
$ws.Range("I86").Value = 2842.9
$ws.Range("J86").Value = 15298.444
$ws.Range("K86").Value = 2842.9
$ws.Range("L86").Value = 15298.444
$ws.Range("M86").Value = -1719.9
$ws.Range("N86").Value = -17544.444
$ws.Range("H89").Value = 8742.895
$ws.Range("I89").Value = 2842.9
$ws.Range("J89").Value = 15298.444
$ws.Range("K89").Value = 14214.5
$ws.Range("L89").Value = 76492.22
$ws.Range("M89").Value = -8598.5
$ws.Range("N89").Value = -87724.22
$ws.Range("H103").Value = 100500790
$ws.Range("I103").Value = 167500160
$ws.Range("K103").Value = 502500480
$ws.Range("M103").Value = -502499894
$ws.Range("H132").Value = 2171.0952
$ws.Range("I132").Value = 2248.15
$ws.Range("J132").Value = 630
$ws.Range("K132").Value = 6744.450000000001
$ws.Range("L132").Value = 1890
$ws.Range("M132").Value = -4214.450000000001
$ws.Range("N132").Value = -6950
$ws.Range("H137").Value = 1477.3549
$ws.Range("I137").Value = 1172.96
$ws.Range("J137").Value = 2745.6667
$ws.Range("K137").Value = 3518.88
$ws.Range("L137").Value = 8237.000100000001
$ws.Range("M137").Value = -968.8800000000001
$ws.Range("N137").Value = -13337.0001
$ws.Range("H141").Value = 2628.5
$ws.Range("I141").Value = 2119.5
$ws.Range("J141").Value = 3901
$ws.Range("K141").Value = 6358.5
$ws.Range("L141").Value = 11703
$ws.Range("M141").Value = -1178.5
$ws.Range("N141").Value = -22063

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3584.8865
$ws.Range("I32").Value = 2853.9768
$ws.Range("K32").Value = 2853.9768
$ws.Range("M32").Value = -2566.9768
$ws.Range("H74").Value = 52632570
$ws.Range("I74").Value = 66667184
$ws.Range("K74").Value = 66667184
$ws.Range("M74").Value = -66666310
$ws.Range("H77").Value = 52632570
$ws.Range("I77").Value = 66667184
$ws.Range("K77").Value = 333335920
$ws.Range("M77").Value = -333331552
$ws.Range("H132").Value = 11000.019
$ws.Range("I132").Value = 1495.6154
$ws.Range("K132").Value = 4486.8462
$ws.Range("M132").Value = -1956.8462

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10972.15
$ws.Range("I31").Value = 17569.9
$ws.Range("J31").Value = 4374.4
$ws.Range("K31").Value = 17569.9
$ws.Range("L31").Value = 4374.4
$ws.Range("M31").Value = -17274.9
$ws.Range("N31").Value = -4964.4
$ws.Range("H34").Value = 10972.15
$ws.Range("I34").Value = 17569.9
$ws.Range("J34").Value = 4374.4
$ws.Range("K34").Value = 17569.9
$ws.Range("L34").Value = 4374.4
$ws.Range("M34").Value = -17367.9
$ws.Range("N34").Value = -4778.4
$ws.Range("H86").Value = 14956.091
$ws.Range("I86").Value = 2807.375
$ws.Range("K86").Value = 2807.375
$ws.Range("M86").Value = -1684.375
$ws.Range("H89").Value = 14956.091
$ws.Range("I89").Value = 2807.375
$ws.Range("K89").Value = 14036.875
$ws.Range("M89").Value = -8420.875
$ws.Range("H99").Value = 5026.522
$ws.Range("I99").Value = 3629.3572
$ws.Range("J99").Value = 7199.8887
$ws.Range("K99").Value = 3629.3572
$ws.Range("L99").Value = 7199.8887
$ws.Range("M99").Value = -2131.3572
$ws.Range("N99").Value = -10195.8887
$ws.Range("H107").Value = 676.9545000000001
$ws.Range("I107").Value = 279.0909
$ws.Range("J107").Value = 1074.8182
$ws.Range("K107").Value = 279.0909
$ws.Range("L107").Value = 1074.8182
$ws.Range("M107").Value = 1640.9091
$ws.Range("N107").Value = -4914.8182
$ws.Range("H126").Value = 5026.522
$ws.Range("I126").Value = 3629.3572
$ws.Range("J126").Value = 7199.8887
$ws.Range("K126").Value = 10888.0716
$ws.Range("L126").Value = 21599.6661
$ws.Range("M126").Value = -8418.071599999999
$ws.Range("N126").Value = -26539.6661

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1044.6842
$ws.Range("I5").Value = 850
$ws.Range("J5").Value = 1186.2727
$ws.Range("K5").Value = 2550
$ws.Range("L5").Value = 3558.8181
$ws.Range("M5").Value = -2438
$ws.Range("N5").Value = -3782.8181
$ws.Range("H131").Value = 780.76
$ws.Range("J131").Value = 780.76
$ws.Range("L131").Value = 2342.28
$ws.Range("N131").Value = -12422.28
$ws.Range("H135").Value = 1044.6842
$ws.Range("I135").Value = 850
$ws.Range("J135").Value = 1186.2727
$ws.Range("K135").Value = 7650
$ws.Range("L135").Value = 10676.4543
$ws.Range("M135").Value = -5115
$ws.Range("N135").Value = -15746.4543

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1535.8572
$ws.Range("I122").Value = 1563.909
$ws.Range("J122").Value = 1433
$ws.Range("K122").Value = 4691.727000000001
$ws.Range("L122").Value = 4299
$ws.Range("M122").Value = -2241.727000000001
$ws.Range("N122").Value = -9199
$ws.Range("H126").Value = 6231.857
$ws.Range("I126").Value = 6193.25
$ws.Range("J126").Value = 6283.3335
$ws.Range("K126").Value = 18579.75
$ws.Range("L126").Value = 18850.0005
$ws.Range("M126").Value = -16109.75
$ws.Range("N126").Value = -23790.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4189.1514
$ws.Range("I7").Value = 3430.1538
$ws.Range("J7").Value = 4682.5
$ws.Range("K7").Value = 3430.1538
$ws.Range("L7").Value = 4682.5
$ws.Range("M7").Value = -3318.1538
$ws.Range("N7").Value = -4906.5
$ws.Range("H82").Value = 2493.3914
$ws.Range("I82").Value = 2288.077
$ws.Range("J82").Value = 2760.3
$ws.Range("K82").Value = 2288.077
$ws.Range("L82").Value = 2760.3
$ws.Range("M82").Value = -1927.077
$ws.Range("N82").Value = -3482.3
$ws.Range("H85").Value = 2493.3914
$ws.Range("I85").Value = 2288.077
$ws.Range("J85").Value = 2760.3
$ws.Range("K85").Value = 2288.077
$ws.Range("L85").Value = 2760.3
$ws.Range("M85").Value = -1040.077
$ws.Range("N85").Value = -5256.3
$ws.Range("H126").Value = 4189.1514
$ws.Range("I126").Value = 3430.1538
$ws.Range("J126").Value = 4682.5
$ws.Range("K126").Value = 10290.4614
$ws.Range("L126").Value = 14047.5
$ws.Range("M126").Value = -7820.4614
$ws.Range("N126").Value = -18987.5
$ws.Range("H132").Value = 6000.6665
$ws.Range("I132").Value = 6004
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 18012
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -15482
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 608.7143
$ws.Range("I107").Value = 656.3
$ws.Range("J107").Value = 489.75
$ws.Range("K107").Value = 1968.9
$ws.Range("L107").Value = 1469.25
$ws.Range("M107").Value = -48.89999999999986
$ws.Range("N107").Value = -5309.25
$ws.Range("H126").Value = 1096.4615
$ws.Range("I126").Value = 851.3333
$ws.Range("J126").Value = 1170
$ws.Range("K126").Value = 2553.9999
$ws.Range("L126").Value = 3510
$ws.Range("M126").Value = -83.9998999999998
$ws.Range("N126").Value = -8450
